$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dFBA net components")

# Current header (A1:H1): Id | Name | dFBA net reaction | Coefficient | Species | Database references | Comments | References
# New header (A1:I1):     Id | Name | dFBA net reaction | Species | Value | Units | Database references | Comments | References

# Shift "dFBA net reaction | Coefficient | Species | Database references | Comments | References"
# one column to the right (C1:H1 -> D1:I1), carrying formatting along, to make room
# for the new "Units" column.
$ws.Range("C1:H1").Copy($ws.Range("D1:I1"))

# Reorder / rename the "Coefficient" and "Species" columns and add "Units".
$ws.Range("D1").Value = "Species"
$ws.Range("E1").Value = "Value"
$ws.Range("F1").Value = "Units"

$ws.Activate()

$wb.Save()
